$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.008.89"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.305.22"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.76"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.52"
$ws.Range("E6").Value = "  -3.47%  "
$ws.Range("E7").Value = "  +3.94%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  +0.62%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("E12").Value = "  -1.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.95"
$ws.Range("E13").Value = "  +0.51%  "
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "2.663.69"
$ws.Range("E15").Value = "  -0.31%  "
$ws.Range("D16").Value = "2.273.57"
$ws.Range("E16").Value = "  -0.93%  "
$ws.Range("E17").Value = "  -3.09%  "
$ws.Range("D18").Value = "42.912.62"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("E19").Value = "  +5.85%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.17"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.23"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.41"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.79"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.14"
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.13"
$ws.Range("E29").Value = "  -1.23%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.03"
$ws.Range("E30").Value = "  -12.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.37"
$ws.Range("E31").Value = "  -3.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.23"
$ws.Range("E32").Value = "  +3.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("E34").Value = "  +1.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "18.22"
$ws.Range("E35").Value = "  +5.08%  "
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E37").Value = "  -0.38%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.102"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.79"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("D42").Value = "1.998.20"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.12"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("E45").Value = "  -7.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.55"
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.82"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "54.61"
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("D49").Value = "2.530.05"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.20"
$ws.Range("E51").Value = "  +4.06%  "
